# Applies numeric value updates to the Odin_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 879.4
$ws.Range("I6").Value = 950
$ws.Range("K6").Value = 2850
$ws.Range("M6").Value = -2738
# Row 8
$ws.Range("H8").Value = 30.875
$ws.Range("I8").Value = 30.875
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 92.625
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 46.375
$ws.Range("N8").ClearContents()
# Row 19
$ws.Range("H19").Value = 2380.4614
$ws.Range("J19").Value = 3357
$ws.Range("L19").Value = 3357
$ws.Range("N19").Value = -3707
# Row 31
$ws.Range("H31").Value = 1399.8889
$ws.Range("I31").Value = 1399.8889
$ws.Range("K31").Value = 4199.6667
$ws.Range("M31").Value = -3969.6667
# Row 38
$ws.Range("H38").Value = 889.7
$ws.Range("J38").Value = 6999
$ws.Range("L38").Value = 20997
$ws.Range("N38").Value = -21741
# Row 39
$ws.Range("H39").Value = 117.3
$ws.Range("I39").Value = 46.5
$ws.Range("J39").Value = 400.5
$ws.Range("K39").Value = 139.5
$ws.Range("L39").Value = 1201.5
$ws.Range("M39").Value = 156.5
$ws.Range("N39").Value = -1793.5
# Row 42
$ws.Range("H42").Value = 2002.6666
$ws.Range("I42").Value = 2002.6666
$ws.Range("K42").Value = 6007.9998
$ws.Range("M42").Value = -5777.9998
# Row 70
$ws.Range("H70").Value = 920.2941
$ws.Range("I70").Value = 645.5833
$ws.Range("J70").Value = 1579.6
$ws.Range("K70").Value = 1936.7499
$ws.Range("L70").Value = 4738.799999999999
$ws.Range("M70").Value = -1666.7499
$ws.Range("N70").Value = -5278.799999999999
# Row 73
$ws.Range("H73").Value = 920.2941
$ws.Range("I73").Value = 645.5833
$ws.Range("J73").Value = 1579.6
$ws.Range("K73").Value = 1936.7499
$ws.Range("L73").Value = 4738.799999999999
$ws.Range("M73").Value = -1000.7499
$ws.Range("N73").Value = -6610.799999999999
# Row 76
$ws.Range("H76").Value = 250010670
$ws.Range("I76").Value = 333342560
$ws.Range("J76").Value = 15000
$ws.Range("K76").Value = 333342560
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -333342245
$ws.Range("N76").Value = -15630
# Row 79
$ws.Range("H79").Value = 250010670
$ws.Range("I79").Value = 333342560
$ws.Range("J79").Value = 15000
$ws.Range("K79").Value = 333342560
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -333341468
$ws.Range("N79").Value = -17184
# Row 88
$ws.Range("H88").Value = 2667.1
$ws.Range("J88").Value = 5166.3335
$ws.Range("L88").Value = 5166.3335
$ws.Range("N88").Value = -5978.3335
# Row 91
$ws.Range("H91").Value = 2667.1
$ws.Range("J91").Value = 5166.3335
$ws.Range("L91").Value = 5166.3335
$ws.Range("N91").Value = -7974.3335
# Row 99
$ws.Range("H99").Value = 100003464
$ws.Range("I99").Value = 888.625
$ws.Range("J99").Value = 500013760
$ws.Range("K99").Value = 2665.875
$ws.Range("L99").Value = 1500041280
$ws.Range("M99").Value = -1167.875
$ws.Range("N99").Value = -1500044276
# Row 138
$ws.Range("H138").Value = 4199.383
$ws.Range("I138").Value = 1344.8422
$ws.Range("K138").Value = 4034.5266
$ws.Range("M138").Value = 1105.4734

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 1740.625
$ws.Range("I63").Value = 954.1667
$ws.Range("K63").Value = 954.1667
$ws.Range("M63").Value = -268.1667
# Row 66
$ws.Range("H66").Value = 1740.625
$ws.Range("I66").Value = 954.1667
$ws.Range("K66").Value = 4770.8335
$ws.Range("M66").Value = -1338.8335
# Row 97
$ws.Range("H97").Value = 1000.2963
$ws.Range("J97").Value = 1687
$ws.Range("L97").Value = 1687
$ws.Range("N97").Value = -2679
# Row 102
$ws.Range("H102").Value = 2253.6365
$ws.Range("I102").Value = 1677.9
$ws.Range("K102").Value = 1677.9
$ws.Range("M102").Value = -55.90000000000009
# Row 132
$ws.Range("H132").Value = 470445.56
$ws.Range("I132").Value = 550317.7
$ws.Range("J132").Value = 91052.836
$ws.Range("K132").Value = 1650953.1
$ws.Range("L132").Value = 273158.508
$ws.Range("M132").Value = -1648423.1
$ws.Range("N132").Value = -278218.508

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1005111.75
$ws.Range("I134").Value = 1195084.5
$ws.Range("J134").Value = 7754.75
$ws.Range("K134").Value = 3585253.5
$ws.Range("L134").Value = 23264.25
$ws.Range("M134").Value = -3582718.5
$ws.Range("N134").Value = -28334.25

$ws = $wb.Worksheets.Item("CRP")
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# Row 105
$ws.Range("H105").Value = 53848644
$ws.Range("J105").Value = 2416.1667
$ws.Range("L105").Value = 2416.1667
$ws.Range("N105").Value = -5910.1667

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1675.8572
$ws.Range("I34").Value = 286.4
$ws.Range("J34").Value = 5149.5
$ws.Range("K34").Value = 859.1999999999999
$ws.Range("L34").Value = 15448.5
$ws.Range("M34").Value = -775.1999999999999
$ws.Range("N34").Value = -15616.5

$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 80
$ws.Range("H80").Value = 166674050
$ws.Range("I80").Value = 333337000
$ws.Range("J80").Value = 11101.667
$ws.Range("K80").Value = 333337000
$ws.Range("L80").Value = 11101.667
$ws.Range("M80").Value = -333336002
$ws.Range("N80").Value = -13097.667
# Row 83
$ws.Range("H83").Value = 166674050
$ws.Range("I83").Value = 333337000
$ws.Range("J83").Value = 11101.667
$ws.Range("K83").Value = 1666685000
$ws.Range("L83").Value = 55508.335
$ws.Range("M83").Value = -1666680008
$ws.Range("N83").Value = -65492.335
# Row 132
$ws.Range("H132").Value = 43483330
$ws.Range("I132").Value = 71433680
$ws.Range("K132").Value = 214301040
$ws.Range("M132").Value = -214298510
# Row 139
$ws.Range("H139").Value = 199999
$ws.Range("J139").Value = 199999
$ws.Range("L139").Value = 199999
$ws.Range("N139").Value = -210279

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1340.4286
$ws.Range("J68").Value = 1349.5
$ws.Range("L68").Value = 1349.5
$ws.Range("N68").Value = -2847.5
# Row 71
$ws.Range("H71").Value = 1340.4286
$ws.Range("J71").Value = 1349.5
$ws.Range("L71").Value = 6747.5
$ws.Range("N71").Value = -14235.5
# Row 82
$ws.Range("H82").Value = 5749.8
$ws.Range("J82").Value = 8283
$ws.Range("L82").Value = 8283
$ws.Range("N82").Value = -9005
# Row 85
$ws.Range("H85").Value = 5749.8
$ws.Range("J85").Value = 8283
$ws.Range("L85").Value = 8283
$ws.Range("N85").Value = -10779
# Row 100
$ws.Range("H100").Value = 3416.4119
$ws.Range("I100").Value = 5829.8335
$ws.Range("J100").Value = 2100
$ws.Range("K100").Value = 5829.8335
$ws.Range("L100").Value = 2100
$ws.Range("M100").Value = -5288.8335
$ws.Range("N100").Value = -3182
# Row 136
$ws.Range("H136").Value = 23812608
$ws.Range("I136").Value = 12198338
$ws.Range("J136").Value = 45457384
$ws.Range("K136").Value = 36595014
$ws.Range("L136").Value = 136372152
$ws.Range("M136").Value = -36592464
$ws.Range("N136").Value = -136377252

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3459.9746
$ws.Range("I132").Value = 2884.1096
$ws.Range("K132").Value = 8652.328799999999
$ws.Range("M132").Value = -6122.328799999999
